# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.975.08'
$ws.Range('E2').Value = '  -0.95%  '

# Row 3: D3, E3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.540.89'
$ws.Range('E3').Value = '  -1.04%  '

# Row 4: E4
$ws.Range('E4').Value = '  -0.09%  '

# Row 5: D5, E5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.69'
$ws.Range('E5').Value = '  +3.20%  '

# Row 6: D6, E6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.86'
$ws.Range('E6').Value = '  -0.75%  '

# Row 7: D7, E7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.536.72'
$ws.Range('E7').Value = '  -0.86%  '

# Row 8: D8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.615'

# Row 9: E9
$ws.Range('E9').Value = '  -0.07%  '

# Row 10: D10, E10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.213'
$ws.Range('E10').Value = '  +5.41%  '

# Row 11: D11, E11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.642'
$ws.Range('E11').Value = '  -1.63%  '

# Row 12: D12, E12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.72'
$ws.Range('E12').Value = '  -2.10%  '

# Row 13: E13
$ws.Range('E13').Value = '  -1.83%  '

# Row 14: E14
$ws.Range('E14').Value = '  -1.61%  '

# Row 15: D15, E15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.104.32'
$ws.Range('E15').Value = '  -1.02%  '

# Row 16: D16, E16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '70.074.77'
$ws.Range('E16').Value = '  -0.86%  '

# Row 17: B17, C17, D17, E17
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '582.24'
$ws.Range('E17').Value = '  +2.16%  '

# Row 18: B18, C18, D18, E18
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.65'
$ws.Range('E18').Value = '  +1.09%  '

# Row 19: B19, C19, D19, E19
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.546.72'
$ws.Range('E19').Value = '  -1.08%  '

# Row 20: B20, C20, D20, E20
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.90'
$ws.Range('E20').Value = '  -2.91%  '

# Row 21: E21
$ws.Range('E21').Value = '  -0.34%  '

# Row 22: D22, E22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.992'
$ws.Range('E22').Value = '  -3.00%  '

# Row 23: D23, E23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.47'
$ws.Range('E23').Value = '  -1.77%  '

# Row 24: D24, E24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.68'
$ws.Range('E24').Value = '  -0.73%  '

# Row 25: D25, E25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.87'
$ws.Range('E25').Value = '  -1.87%  '

# Row 26: D26, E26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '95.75'
$ws.Range('E26').Value = '  -0.47%  '

# Row 27: D27, E27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.99'
$ws.Range('E27').Value = '  -0.43%  '

# Row 28: D28, E28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.98'
$ws.Range('E28').Value = '  -4.78%  '

# Row 29: D29, E29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.47'
$ws.Range('E29').Value = '  +3.04%  '

# Row 30: D30, E30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.08'
$ws.Range('E30').Value = '  -1.02%  '

# Row 31: D31, E31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.02'
$ws.Range('E31').Value = '  -5.03%  '

# Row 32: D32, E32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.16'
$ws.Range('E32').Value = '  -3.11%  '

# Row 33: E33
$ws.Range('E33').Value = '  -1.45%  '

# Row 34: D34, E34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.10'
$ws.Range('E34').Value = '  -3.13%  '

# Row 35: D35, E35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.63'
$ws.Range('E35').Value = '  +17.46%  '

# Row 36: D36, E36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.25'
$ws.Range('E36').Value = '  -0.36%  '

# Row 37: D37, E37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '536.82'
$ws.Range('E37').Value = '  -5.15%  '

# Row 38: E38
$ws.Range('E38').Value = '  -3.55%  '

# Row 39: E39
$ws.Range('E39').Value = '  +0.16%  '

# Row 40: D40, E40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.16'
$ws.Range('E40').Value = '  -2.92%  '

# Row 41: E41
$ws.Range('E41').Value = '  +0.14%  '

# Row 42: D42, E42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.523.66'
$ws.Range('E42').Value = '  +5.05%  '

# Row 43: E43
$ws.Range('E43').Value = '  +4.06%  '

# Row 44: E44
$ws.Range('E44').Value = '  +0.64%  '

# Row 45: D45, E45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0455'
$ws.Range('E45').Value = '  +1.85%  '

# Row 46: E46
$ws.Range('E46').Value = '  -2.13%  '

# Row 47: D47, E47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.38'
$ws.Range('E47').Value = '  -5.32%  '

# Row 48: D48, E48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.141'
$ws.Range('E48').Value = '  +2.53%  '

# Row 49: D49, E49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.13'
$ws.Range('E49').Value = '  -3.54%  '

# Row 50: E50
$ws.Range('E50').Value = '  +0.14%  '

# Row 51: E51
$ws.Range('E51').Value = '  -2.81%  '
